$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(16, 1).Value = "Pakistan"
$ws.Cells.Item(16, 2).Value = 160118
$ws.Cells.Item(16, 3).Value = 5358
$ws.Cells.Item(16, 4).Value = 59215
$ws.Cells.Item(16, 5).Value = 97810
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 118
$ws.Cells.Item(16, 8).Value = 3093

$ws.Cells.Item(17, 1).Value = "Mexico"
$ws.Cells.Item(17, 2).Value = 159793
$ws.Cells.Item(17, 3).Value = 4930
$ws.Cells.Item(17, 4).Value = 119355
$ws.Cells.Item(17, 5).Value = 21358
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 770
$ws.Cells.Item(17, 8).Value = 19080

$ws.Cells.Item(18, 1).Value = "Francia"
$ws.Cells.Item(18, 2).Value = 158174
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 73667
$ws.Cells.Item(18, 5).Value = 54932
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 29575

$ws.Cells.Item(56, 1).Value = "Kazajistan"
$ws.Cells.Item(56, 2).Value = 15877
$ws.Cells.Item(56, 3).Value = 335
$ws.Cells.Item(56, 4).Value = 9920
$ws.Cells.Item(56, 5).Value = 5860
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 97

$ws.Cells.Item(80, 1).Value = "Haiti"
$ws.Cells.Item(80, 2).Value = 4688
$ws.Cells.Item(80, 3).Value = 141
$ws.Cells.Item(80, 4).Value = 24
$ws.Cells.Item(80, 5).Value = 4582
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 2
$ws.Cells.Item(80, 8).Value = 82

$ws.Cells.Item(81, 1).Value = "Guinea"
$ws.Cells.Item(81, 2).Value = 4668
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 4).Value = 3364
$ws.Cells.Item(81, 5).Value = 1278
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 26

$ws.Cells.Item(96, 1).Value = "Kirguistan"
$ws.Cells.Item(96, 2).Value = 2657
$ws.Cells.Item(96, 3).Value = 95
$ws.Cells.Item(96, 4).Value = 1933
$ws.Cells.Item(96, 5).Value = 693
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 31

$ws.Cells.Item(97, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(97, 2).Value = 2564
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = 402
$ws.Cells.Item(97, 5).Value = 2144
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 18

$ws.Cells.Item(163, 1).Value = "Mongolia"
$ws.Cells.Item(163, 2).Value = 201
$ws.Cells.Item(163, 3).Value = 4
$ws.Cells.Item(163, 4).Value = 127
$ws.Cells.Item(163, 5).Value = 74
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 0

$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0

$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 06:22"
